# Bill Summary sheet update:
# A new line item ("Short point (up to 3 mtr.)") is inserted as row 9,
# pushing all subsequent rows down by one. Quantities/amounts for the
# line items and the grand-total rows are recalculated accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Insert a new row at position 9 (shifts old rows 9-20 to 10-21)
# ---------------------------------------------------------------
$ws.Rows.Item(9).Insert()

# ---------------------------------------------------------------
# 2. Row 8 - "Qty executed upto date" changes
# ---------------------------------------------------------------
$ws.Range("C8").Value = 87

# ---------------------------------------------------------------
# 3. New row 9 - "Short point (up to 3 mtr.)" line item
# ---------------------------------------------------------------
$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = "P. point"
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 82
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "Short point (up to 3 mtr.)"
$ws.Range("F9").Value = 256
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "20992.00"
$ws.Range("H9").Value = 0

# ---------------------------------------------------------------
# 4. Rows 10-17 (previously 9-16, shifted down by the insert) keep
#    their original text but get updated quantities / amounts
# ---------------------------------------------------------------
# Row 10 - Medium point (up to 6 mtr.)
$ws.Range("C10").Value = 61
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "28792.00"

# Row 11 - Long point (up to 10 mtr.)
$ws.Range("C11").Value = 6
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "3972.00"

# Row 12 - Rewiring of 3/5 pin 6 amp light plug point
$ws.Range("C12").Value = 98

# Row 13 - On board
$ws.Range("C13").Value = 56
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "7616.00"

# Row 14 - P & F ISI marked switch
$ws.Range("C14").Value = 89
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "2047.00"

# Row 15 - Total
$ws.Range("C15").Value = 87

# Row 16 - Add Tender Premium
$ws.Range("C16").Value = 11

# Row 17 - Grand Total
$ws.Range("C17").Value = 53

# ---------------------------------------------------------------
# 5. Summary rows (previously 18-20, now 19-21) - updated totals
# ---------------------------------------------------------------
# Row 19 - Grand Total Rs.
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "63419.00"
$ws.Range("H19").NumberFormat = "@"
$ws.Range("H19").Value = "63419.00"

# Row 21 - NET PAYABLE AMOUNT Rs.
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "63419.00"
$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value = "63419.00"
